$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New data for the transaction list table (rows 2-11, columns A-C)
$data = @(
    @(0, "Connexion", "00;Bob;sonMotDePasse"),
    @(1, "Connexion réussie", "01;"),
    @(2, "Erreur lors de l'autentification ", "02;"),
    @(3, "Carte piochée", "03;objet carte"),
    @(4, "Carte jouée", "04;objet carte; nom du joueur; nouveau total"),
    @(5, "Carte choisie", "To DO"),
    @(6, "Partie terminée", "05;V (ou D)"),
    @(7, "Obtenir le leaderboard", "06;player1;999;player2;99;…;player10;1"),
    @(8, "Carte non-valide", "07;"),
    @(9, "Dire à un joueur que c'est à son tour.", "08;")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Resize the table / list object to include the newly added row
$table = $ws.ListObjects.Item("Tableau1")
$table.Resize($ws.Range("A1:C11"))
